# Adding the changes we made on may 9th
#
# accelerometer_selected.xlsx originally held 20 samples (rows 2-21).
# This change prepends 7 new samples (shifting the existing 20 rows
# down to rows 9-28) and appends 3 more samples at the end (rows 29-31),
# for a new total of 30 samples (rows 2-31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the existing data (rows 2-21, columns A:C) before moving anything,
# so we don't clobber values while shifting them down.
$existing = @()
for ($r = 2; $r -le 21; $r++) {
    $existing += ,@($ws.Cells.Item($r, 1).Value2, $ws.Cells.Item($r, 2).Value2, $ws.Cells.Item($r, 3).Value2)
}

# Re-write the existing 20 rows shifted down by 7 (new rows 9-28).
for ($i = 0; $i -lt $existing.Count; $i++) {
    $rowNum = 9 + $i
    $ws.Cells.Item($rowNum, 1).Value = $existing[$i][0]
    $ws.Cells.Item($rowNum, 2).Value = $existing[$i][1]
    $ws.Cells.Item($rowNum, 3).Value = $existing[$i][2]
}

# New data block inserted at the top (new rows 2-8).
$topInsert = @(
    @(3.73009729385376, -2.067999362945557, 2.923094749450684),
    @(1.413437724113464, -7.171976566314697, 6.755977630615234),
    @(-4.193170547485352, -4.375148296356201, 1.568653106689453),
    @(-1.307081580162048, -2.807691097259521, -1.606552600860596),
    @(-1.827142477035522, 0.1487736701965332, 4.701539993286133),
    @(-1.866428852081299, 4.58729076385498, 1.570873260498047),
    @(1.589986324310303, -0.6990594863891602, 7.147370338439941)
)

for ($i = 0; $i -lt $topInsert.Count; $i++) {
    $rowNum = 2 + $i
    $ws.Cells.Item($rowNum, 1).Value = $topInsert[$i][0]
    $ws.Cells.Item($rowNum, 2).Value = $topInsert[$i][1]
    $ws.Cells.Item($rowNum, 3).Value = $topInsert[$i][2]
}

# New data block appended at the end (rows 29-31).
$bottomAppend = @(
    @(-12.98141479492188, -5.830618858337402, -7.131386756896973),
    @(15.02213287353516, -47.13114547729492, 8.131996154785156),
    @(-2.31139087677002, 14.0579719543457, 5.640069961547852)
)

for ($i = 0; $i -lt $bottomAppend.Count; $i++) {
    $rowNum = 29 + $i
    $ws.Cells.Item($rowNum, 1).Value = $bottomAppend[$i][0]
    $ws.Cells.Item($rowNum, 2).Value = $bottomAppend[$i][1]
    $ws.Cells.Item($rowNum, 3).Value = $bottomAppend[$i][2]
}
